# Insert a new price record as row 89 on the "Alcachofa" sheet, pushing the
# existing rows 89-141 down to 90-142 (dimension grows from R141 to R142).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 89..141 down by one, duplicating row 89's formatting (date style
# on column D) onto the freshly inserted row.
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with the new observation.
$ws.Cells.Item(89, 1).Value2  = 10
$ws.Cells.Item(89, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(89, 3).Value2  = "La Araucanía"
$ws.Cells.Item(89, 4).Value2  = 44489
$ws.Cells.Item(89, 5).Value2  = 9
$ws.Cells.Item(89, 6).Value2  = 100112013
$ws.Cells.Item(89, 7).Value2  = "Alcachofa"
$ws.Cells.Item(89, 8).Value2  = "Española"
$ws.Cells.Item(89, 9).Value2  = "Primera"
$ws.Cells.Item(89, 10).Value2 = 55
$ws.Cells.Item(89, 11).Value2 = 12000
$ws.Cells.Item(89, 12).Value2 = 12000
$ws.Cells.Item(89, 13).Value2 = 12000
$ws.Cells.Item(89, 14).Value2 = "`$/caja 30 unidades"
$ws.Cells.Item(89, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(89, 16).Value2 = 400
$ws.Cells.Item(89, 17).Value2 = 30
$ws.Cells.Item(89, 18).Value2 = "Hortaliza"
